$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 30.87085333333333
$ws.Range("H2").Value = 92.61256
$ws.Range("I2").Value = 0.2985789950947061
$ws.Range("J2").Value = 0.2985789950947061
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.06447966666666667
$ws.Range("N2").Value = 0.193439
$ws.Range("O2").Value = 0.001101138907643723
$ws.Range("P2").Value = 0.001101138907643722
$ws.Range("Q2").Value = 1.990542332648889
$ws.Range("R2").Value = 17.91488099384
$ws.Range("S2").Value = 0.0003287769485039451
$ws.Range("T2").Value = 0.0003287769485039451

$ws.Range("G3").Value = 30.87085333333333
$ws.Range("H3").Value = 92.61256
$ws.Range("I3").Value = 0.2985789950947061
$ws.Range("J3").Value = 0.2985789950947061
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.3851286666666667
$ws.Range("N3").Value = 1.155386
$ws.Range("O3").Value = 0.00657695954769643
$ws.Range("P3").Value = 0.006576959547696431
$ws.Range("Q3").Value = 11.88925058312889
$ws.Range("R3").Value = 107.00325524816
$ws.Range("S3").Value = 0.001963741972529733
$ws.Range("T3").Value = 0.001963741972529734

$ws.Range("G4").Value = 30.87085333333333
$ws.Range("H4").Value = 92.61256
$ws.Range("I4").Value = 0.2985789950947061
$ws.Range("J4").Value = 0.2985789950947061
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.008175
$ws.Range("N4").Value = 0.024525
$ws.Range("O4").Value = 0.0001396069650378791
$ws.Range("P4").Value = 0.0001396069650378791
$ws.Range("Q4").Value = 0.252369226
$ws.Range("R4").Value = 2.271323034
$ws.Range("S4").Value = 0.0000416837073292317
$ws.Range("T4").Value = 0.00004168370732923171

$ws.Range("G5").Value = 30.87085333333333
$ws.Range("H5").Value = 92.61256
$ws.Range("I5").Value = 0.2985789950947061
$ws.Range("J5").Value = 0.2985789950947061
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 58.099467
$ws.Range("N5").Value = 174.298401
$ws.Range("O5").Value = 0.992182294579622
$ws.Range("P5").Value = 0.992182294579622
$ws.Range("Q5").Value = 1793.58012450184
$ws.Range("R5").Value = 16142.22112051656
$ws.Range("S5").Value = 0.2962447924663432
$ws.Range("T5").Value = 0.2962447924663433

$ws.Range("G6").Value = 33.793597
$ws.Range("H6").Value = 101.380791
$ws.Range("I6").Value = 0.3268474027571036
$ws.Range("J6").Value = 0.3268474027571037
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.06447966666666667
$ws.Range("N6").Value = 0.193439
$ws.Range("O6").Value = 0.001101138907643723
$ws.Range("P6").Value = 0.001101138907643722
$ws.Range("Q6").Value = 2.178999870027667
$ws.Range("R6").Value = 19.610998830249
$ws.Range("S6").Value = 0.0003599043920381449
$ws.Range("T6").Value = 0.0003599043920381449

$ws.Range("G7").Value = 33.793597
$ws.Range("H7").Value = 101.380791
$ws.Range("I7").Value = 0.3268474027571036
$ws.Range("J7").Value = 0.3268474027571037
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.3851286666666667
$ws.Range("N7").Value = 1.155386
$ws.Range("O7").Value = 0.00657695954769643
$ws.Range("P7").Value = 0.006576959547696431
$ws.Range("Q7").Value = 13.01488295448067
$ws.Range("R7").Value = 117.133946590326
$ws.Range("S7").Value = 0.002149662146203113
$ws.Range("T7").Value = 0.002149662146203114

$ws.Range("G8").Value = 33.793597
$ws.Range("H8").Value = 101.380791
$ws.Range("I8").Value = 0.3268474027571036
$ws.Range("J8").Value = 0.3268474027571037
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.008175
$ws.Range("N8").Value = 0.024525
$ws.Range("O8").Value = 0.0001396069650378791
$ws.Range("P8").Value = 0.0001396069650378791
$ws.Range("Q8").Value = 0.276262655475
$ws.Range("R8").Value = 2.486363899275
$ws.Range("S8").Value = 0.00004563017392943254
$ws.Range("T8").Value = 0.00004563017392943255

$ws.Range("G9").Value = 33.793597
$ws.Range("H9").Value = 101.380791
$ws.Range("I9").Value = 0.3268474027571036
$ws.Range("J9").Value = 0.3268474027571037
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 58.099467
$ws.Range("N9").Value = 174.298401
$ws.Range("O9").Value = 0.992182294579622
$ws.Range("P9").Value = 0.992182294579622
$ws.Range("Q9").Value = 1963.389973712799
$ws.Range("R9").Value = 17670.50976341519
$ws.Range("S9").Value = 0.3242922060449329
$ws.Range("T9").Value = 0.324292206044933

$ws.Range("G10").Value = 2.981185666666667
$ws.Range("H10").Value = 8.943557
$ws.Range("I10").Value = 0.02883365130639111
$ws.Range("J10").Value = 0.02883365130639111
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.06447966666666667
$ws.Range("N10").Value = 0.193439
$ws.Range("O10").Value = 0.001101138907643723
$ws.Range("P10").Value = 0.001101138907643722
$ws.Range("Q10").Value = 0.1922258580581111
$ws.Range("R10").Value = 1.730032722523
$ws.Range("S10").Value = 0.00003174985530289949
$ws.Range("T10").Value = 0.00003174985530289949

$ws.Range("G11").Value = 2.981185666666667
$ws.Range("H11").Value = 8.943557
$ws.Range("I11").Value = 0.02883365130639111
$ws.Range("J11").Value = 0.02883365130639111
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.3851286666666667
$ws.Range("N11").Value = 1.155386
$ws.Range("O11").Value = 0.00657695954769643
$ws.Range("P11").Value = 0.006576959547696431
$ws.Range("Q11").Value = 1.148140060889111
$ws.Range("R11").Value = 10.333260548002
$ws.Range("S11").Value = 0.0001896377582545186
$ws.Range("T11").Value = 0.0001896377582545186

$ws.Range("G12").Value = 2.981185666666667
$ws.Range("H12").Value = 8.943557
$ws.Range("I12").Value = 0.02883365130639111
$ws.Range("J12").Value = 0.02883365130639111
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.008175
$ws.Range("N12").Value = 0.024525
$ws.Range("O12").Value = 0.0001396069650378791
$ws.Range("P12").Value = 0.0001396069650378791
$ws.Range("Q12").Value = 0.024371192825
$ws.Range("R12").Value = 0.219340735425
$ws.Range("S12").Value = 0.000004025378549845739
$ws.Range("T12").Value = 0.000004025378549845739

$ws.Range("G13").Value = 2.981185666666667
$ws.Range("H13").Value = 8.943557
$ws.Range("I13").Value = 0.02883365130639111
$ws.Range("J13").Value = 0.02883365130639111
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 58.099467
$ws.Range("N13").Value = 174.298401
$ws.Range("O13").Value = 0.992182294579622
$ws.Range("P13").Value = 0.992182294579622
$ws.Range("Q13").Value = 173.205298261373
$ws.Range("R13").Value = 1558.847684352357
$ws.Range("S13").Value = 0.02860823831428384
$ws.Range("T13").Value = 0.02860823831428384

$ws.Range("G14").Value = 35.74694633333333
$ws.Range("H14").Value = 107.240839
$ws.Range("I14").Value = 0.3457399508417991
$ws.Range("J14").Value = 0.3457399508417991
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.06447966666666667
$ws.Range("N14").Value = 0.193439
$ws.Range("O14").Value = 0.001101138907643723
$ws.Range("P14").Value = 0.001101138907643722
$ws.Range("Q14").Value = 2.304951183924556
$ws.Range("R14").Value = 20.744560655321
$ws.Range("S14").Value = 0.000380707711798733
$ws.Range("T14").Value = 0.0003807077117987329

$ws.Range("G15").Value = 35.74694633333333
$ws.Range("H15").Value = 107.240839
$ws.Range("I15").Value = 0.3457399508417991
$ws.Range("J15").Value = 0.3457399508417991
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = 0.6666666666666666
$ws.Range("M15").Value = 0.3851286666666667
$ws.Range("N15").Value = 1.155386
$ws.Range("O15").Value = 0.00657695954769643
$ws.Range("P15").Value = 0.006576959547696431
$ws.Range("Q15").Value = 13.76717377876156
$ws.Range("R15").Value = 123.904564008854
$ws.Range("S15").Value = 0.002273917670709065
$ws.Range("T15").Value = 0.002273917670709065

$ws.Range("G16").Value = 35.74694633333333
$ws.Range("H16").Value = 107.240839
$ws.Range("I16").Value = 0.3457399508417991
$ws.Range("J16").Value = 0.3457399508417991
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.008175
$ws.Range("N16").Value = 0.024525
$ws.Range("O16").Value = 0.0001396069650378791
$ws.Range("P16").Value = 0.0001396069650378791
$ws.Range("Q16").Value = 0.292231286275
$ws.Range("R16").Value = 2.630081576475
$ws.Range("S16").Value = 0.00004826770522936908
$ws.Range("T16").Value = 0.00004826770522936908

$ws.Range("G17").Value = 35.74694633333333
$ws.Range("H17").Value = 107.240839
$ws.Range("I17").Value = 0.3457399508417991
$ws.Range("J17").Value = 0.3457399508417991
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 58.099467
$ws.Range("N17").Value = 174.298401
$ws.Range("O17").Value = 0.992182294579622
$ws.Range("P17").Value = 0.992182294579622
$ws.Range("Q17").Value = 2076.878528844271
$ws.Range("R17").Value = 18691.90675959844
$ws.Range("S17").Value = 0.343037057754062
$ws.Range("T17").Value = 0.343037057754062

Write-Host "done"
